$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Color constants (VBA-style BGR integers matching the workbook's existing fill palette)
$colYellow = 65535    # FFFF00
$colOrange = 49407    # FFC000
$xlLeft = -4131

# --- 1. Row 7: Pleiades line gains "Refactor C#, " ---
$ws.Range("A7").Value = "Pleiades (Refactor C#, Product Management, Git Tools, MVC)"

# --- 2. Old row 19 ("Get More T-Shirts...") gets reworded, still row 19 at this point ---
$ws.Range("A19").Value = "Get More T-Shirts - see links + http://www.aeropostale.com"
$ws.Range("A19").Interior.Color = $colYellow

# --- 3. Old row 20 ("Moby, Bon Iver, ") gets extended, still row 20 at this point ---
$ws.Range("A20").Value = "Moby, Bon Iver, Chicago Mix Tape, etc."

# --- 4. Row 13: budget line gets an addendum ---
$ws.Range("A13").Value = "PLEASE COME UP WITH A BUDGET! & LOG INTO PAYROLL WEBSITE"

# --- 5. Rows 8-10 (Professional section) reshuffle ---
# "Buy Resharper" moves up into row 8 (now ON QUEUE) and loses its bold weight;
# "Get a Monitor" moves down into row 10 (now TODO, alignment reset to general).
$ws.Range("A8").Value = "Buy Resharper - for productivity"
$ws.Range("A8").Font.Bold = $false
$ws.Range("A8").Interior.Color = $colOrange

$ws.Range("B8").Value = "ON QUEUE"
$ws.Range("B8").Font.Bold = $false
$ws.Range("B8").Interior.Color = $colOrange
$ws.Range("B8").HorizontalAlignment = $xlLeft

$ws.Range("A10").Value = "Get a Monitor / De-commission / Replace Big Blue"
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)   # xlPasteFormats -- matches the plain TODO look (no left align)
$excel.CutCopyMode = $false

# --- 6. Insert a new row after "RETURN GI" (row 18) for "Goto Fleetfeet" ---
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "Goto Fleetfeet"
$ws.Range("B19").Value = "ON QUEUE"
